$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before column A, shifting the existing
# Project Name / Emp Name / ... headers from A:I to B:J.
$ws.Columns("A:A").Insert()

# Copy the header formatting (bold white text on blue fill, centered,
# bordered) from the now-shifted "Project Name" header in B1 onto the
# new A1 cell, then set its text.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("A1").Value = "S.No"

# Match the saved selection/view state.
$ws.Range("E2").Select()
